$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dhh"
$ws.Range("C2").Value = "Ptch1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.179395
$ws.Range("H2").Value = 6.538185
$ws.Range("I2").Value = 0.4845018986408914
$ws.Range("J2").Value = 0.4845018986408914
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.822797666666667
$ws.Range("N2").Value = 11.468393
$ws.Range("O2").Value = 0.170387229893941
$ws.Range("P2").Value = 0.170387229893941
$ws.Range("Q2").Value = 8.331386120745
$ws.Range("R2").Value = 74.98247508670501
$ws.Range("S2").Value = 0.08255293638777644
$ws.Range("T2").Value = 0.08255293638777647

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dhh"
$ws.Range("C3").Value = "Ptch1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.179395
$ws.Range("H3").Value = 6.538185
$ws.Range("I3").Value = 0.4845018986408914
$ws.Range("J3").Value = 0.4845018986408914
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.199585666666667
$ws.Range("N3").Value = 21.598757
$ws.Range("O3").Value = 0.3208952095016597
$ws.Range("P3").Value = 0.3208952095016597
$ws.Range("Q3").Value = 15.690741004005
$ws.Range("R3").Value = 141.216669036045
$ws.Range("S3").Value = 0.1554743382683207
$ws.Range("T3").Value = 0.1554743382683207

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dhh"
$ws.Range("C4").Value = "Ptch1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 2.179395
$ws.Range("H4").Value = 6.538185
$ws.Range("I4").Value = 0.4845018986408914
$ws.Range("J4").Value = 0.4845018986408914
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.41355666666667
$ws.Range("N4").Value = 34.24067
$ws.Range("O4").Value = 0.5087175606043992
$ws.Range("P4").Value = 0.5087175606043994
$ws.Range("Q4").Value = 24.87464833155
$ws.Range("R4").Value = 223.87183498395
$ws.Range("S4").Value = 0.2464746239847941
$ws.Range("T4").Value = 0.2464746239847942

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Dhh"
$ws.Range("C5").Value = "Ptch1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7348883333333333
$ws.Range("H5").Value = 2.204665
$ws.Range("I5").Value = 0.1633732264179005
$ws.Range("J5").Value = 0.1633732264179005
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.822797666666667
$ws.Range("N5").Value = 11.468393
$ws.Range("O5").Value = 0.170387229893941
$ws.Range("P5").Value = 0.170387229893941
$ws.Range("Q5").Value = 2.809329405927222
$ws.Range("R5").Value = 25.283964653345
$ws.Range("S5").Value = 0.02783671148818168
$ws.Range("T5").Value = 0.02783671148818169

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dhh"
$ws.Range("C6").Value = "Ptch1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7348883333333333
$ws.Range("H6").Value = 2.204665
$ws.Range("I6").Value = 0.1633732264179005
$ws.Range("J6").Value = 0.1633732264179005
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.199585666666667
$ws.Range("N6").Value = 21.598757
$ws.Range("O6").Value = 0.3208952095016597
$ws.Range("P6").Value = 0.3208952095016597
$ws.Range("Q6").Value = 5.290891511267222
$ws.Range("R6").Value = 47.61802360140499
$ws.Range("S6").Value = 0.05242568571833425
$ws.Range("T6").Value = 0.05242568571833426

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dhh"
$ws.Range("C7").Value = "Ptch1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7348883333333333
$ws.Range("H7").Value = 2.204665
$ws.Range("I7").Value = 0.1633732264179005
$ws.Range("J7").Value = 0.1633732264179005
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.41355666666667
$ws.Range("N7").Value = 34.24067
$ws.Range("O7").Value = 0.5087175606043992
$ws.Range("P7").Value = 0.5087175606043994
$ws.Range("Q7").Value = 8.387689636172221
$ws.Range("R7").Value = 75.48920672555001
$ws.Range("S7").Value = 0.08311082921138452
$ws.Range("T7").Value = 0.08311082921138456

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Dhh"
$ws.Range("C8").Value = "Ptch1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.583934333333333
$ws.Range("H8").Value = 4.751803
$ws.Range("I8").Value = 0.3521248749412083
$ws.Range("J8").Value = 0.3521248749412082
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.822797666666667
$ws.Range("N8").Value = 11.468393
$ws.Range("O8").Value = 0.170387229893941
$ws.Range("P8").Value = 0.170387229893941
$ws.Range("Q8").Value = 6.055060473619889
$ws.Range("R8").Value = 54.495544262579
$ws.Range("S8").Value = 0.05999758201798287
$ws.Range("T8").Value = 0.05999758201798288

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Dhh"
$ws.Range("C9").Value = "Ptch1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.583934333333333
$ws.Range("H9").Value = 4.751803
$ws.Range("I9").Value = 0.3521248749412083
$ws.Range("J9").Value = 0.3521248749412082
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.199585666666667
$ws.Range("N9").Value = 21.598757
$ws.Range("O9").Value = 0.3208952095016597
$ws.Range("P9").Value = 0.3208952095016597
$ws.Range("Q9").Value = 11.40367092320789
$ws.Range("R9").Value = 102.633038308871
$ws.Range("S9").Value = 0.1129951855150047
$ws.Range("T9").Value = 0.1129951855150047

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Dhh"
$ws.Range("C10").Value = "Ptch1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.583934333333333
$ws.Range("H10").Value = 4.751803
$ws.Range("I10").Value = 0.3521248749412083
$ws.Range("J10").Value = 0.3521248749412082
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 11.41355666666667
$ws.Range("N10").Value = 34.24067
$ws.Range("O10").Value = 0.5087175606043992
$ws.Range("P10").Value = 0.5087175606043994
$ws.Range("Q10").Value = 18.07832426977889
$ws.Range("R10").Value = 162.70491842801
$ws.Range("S10").Value = 0.1791321074082206
$ws.Range("T10").Value = 0.1791321074082207
